$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1: bump the "last status check" timestamp from 01:00 to 01:15
$ws.Range("F1").Value = "Last status check on: 27.01.2022 01:15"

# Row 5 ("Makro"): the current/old prices swapped, the delta flipped sign
# and is now recorded as text, and the "old date" cell switched from a
# numeric Excel datetime to a plain text timestamp string.
$ws.Range("B5").Value = 34.5
$ws.Range("C5").Value = 34.9

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "-0.4"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2022-01-27 01:15:08"
$ws.Range("E5").Style = "Normal"
